# Insert two new respondent rows into the "base_nettoyee" sheet, just above
# the row currently holding hkboboroto@gmail.com (row 191), pushing the
# existing rows down by two. Then populate the two new rows with the new
# respondents' data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 2 blank rows at position 191 (rows 191 & 192), shifting the
# pre-existing data (previously starting at row 191) down to rows 193+.
$ws.Range("A191:A192").EntireRow.Insert()

# --- New row 191: pierreyvon7@gmail.com ---
$ws.Range("A191").Value = "pierreyvon7@gmail.com"
$ws.Range("B191").Value = 0
$ws.Range("C191").Value = "1 3 5 10 14"
$ws.Range("E191").Value = "commerce/gestion/economie"
$ws.Range("F191").Value = "Homme"
$ws.Range("G191").Value = "25-35 ans"
$ws.Range("H191").Value = "Togo"
$ws.Range("I191").Value = 1
$ws.Range("J191").Value = 0
$ws.Range("K191").Value = 1
$ws.Range("L191").Value = 0
$ws.Range("M191").Value = 1
$ws.Range("N191").Value = 0
$ws.Range("O191").Value = 0
$ws.Range("P191").Value = 0
$ws.Range("Q191").Value = 0
$ws.Range("R191").Value = 1
$ws.Range("S191").Value = 0
$ws.Range("T191").Value = 0
$ws.Range("U191").Value = 0
$ws.Range("V191").Value = 1
$ws.Range("W191").Value = "commerce, economie et comptabilite"
$ws.Range("X191").Value = 5

# --- New row 192: anamingcrepin5@gmail.com ---
$ws.Range("A192").Value = "anamingcrepin5@gmail.com"
$ws.Range("B192").Value = 0
$ws.Range("C192").Value = "5 8 12 13 15"
$ws.Range("E192").Value = "commerce/gestion/economie"
$ws.Range("F192").Value = "Homme"
$ws.Range("G192").Value = "25-35 ans"
$ws.Range("H192").Value = "Togo"
$ws.Range("I192").Value = 0
$ws.Range("J192").Value = 0
$ws.Range("K192").Value = 0
$ws.Range("L192").Value = 0
$ws.Range("M192").Value = 1
$ws.Range("N192").Value = 0
$ws.Range("O192").Value = 0
$ws.Range("P192").Value = 1
$ws.Range("Q192").Value = 0
$ws.Range("R192").Value = 0
$ws.Range("S192").Value = 0
$ws.Range("T192").Value = 1
$ws.Range("U192").Value = 1
$ws.Range("V192").Value = 0
$ws.Range("W192").Value = "commerce, economie et comptabilite"
$ws.Range("X192").Value = 4
